$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.022559779323172841
$ws.Range("B1").Value = -0.023196238584140032
$ws.Range("A2").Value = -0.037900829027244683
$ws.Range("B2").Value = -0.025623603773802562
$ws.Range("A3").Value = -0.013782967777531232
$ws.Range("B3").Value = -0.013281671387067855
$ws.Range("A4").Value = -0.061818684486596109
$ws.Range("B4").Value = -0.061754065911548597
$ws.Range("A5").Value = -0.022366786310500587
$ws.Range("B5").Value = -0.022282449892954616
